$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 278 (shifts existing rows 278..319 down to 279..320)
$ws.Rows(278).Insert()

# Populate the newly inserted row 278 with the weekly record
# (same market/category template as the neighboring rows, new date + prices)
$ws.Range("A278").Value = 5
$ws.Range("B278").Value = "Macroferia Regional de Talca"
$ws.Range("C278").Value = "Maule"
$ws.Range("D278").Value = 44984
$ws.Range("E278").Value = 7
$ws.Range("F278").Value = 100112024
$ws.Range("G278").Value = "Choclo"
$ws.Range("H278").Value = "Choclero"
$ws.Range("I278").Value = "Primera"
$ws.Range("J278").Value = 40000
$ws.Range("K278").Value = 400
$ws.Range("L278").Value = 450
$ws.Range("M278").Value = 425
$ws.Range("N278").Value = '$/unidad'
$ws.Range("O278").Value = "Región del Maule"
$ws.Range("P278").Value = 425
$ws.Range("Q278").Value = 1
$ws.Range("R278").Value = "Hortaliza"
